# Passenger rail maintenance costs: divide by 10 to align with
# multiplying vehicle counts by 10 (see "Cost Data" sheet, row 88).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cost Data")

# Update the two formulas that feed the passenger-rail maintenance row.
$ws.Range("B88").Formula = "=B54/10"
$ws.Range("C88").Formula = "=B55/10"

# These cells picked up a distinct style in the source workbook: same
# green fill as before, but number format switched from a whole-dollar
# format to a two-decimal dollar format.
$ws.Range("B88:C88").NumberFormat = """$""#,##0.00"
$ws.Range("B88:C88").Interior.Color = 5296274

$wb.Save()
